$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = '93.044.93'; E = '  +1.51%  ' },
    @{ Row = 3; D = '3.094.42'; E = '  -1.00%  ' },
    @{ Row = 4; E = '  +0.08%  ' },
    @{ Row = 5; D = '235.46'; E = '  -4.34%  ' },
    @{ Row = 6; D = '611.01'; E = '  -1.05%  ' },
    @{ Row = 7; D = '1.13'; E = '  +2.23%  ' },
    @{ Row = 8; D = '0.385'; E = '  +0.48%  ' },
    @{ Row = 9; E = '  +0.00%  ' },
    @{ Row = 10; D = '0.822'; E = '  +11.53%  ' },
    @{ Row = 11; D = '3.091.66'; E = '  -0.90%  ' },
    @{ Row = 12; E = '  -3.10%  ' },
    @{ Row = 13; D = '0.0000242'; E = '  -3.46%  ' },
    @{ Row = 14; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '92.975.89'; E = '  +1.60%  ' },
    @{ Row = 15; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '34.62'; E = '  -0.69%  ' },
    @{ Row = 16; D = '5.39'; E = '  -3.79%  ' },
    @{ Row = 17; D = '3.678.86'; E = '  -0.63%  ' },
    @{ Row = 18; D = '3.105.01'; E = '  -1.29%  ' },
    @{ Row = 19; D = '3.68'; E = '  -0.89%  ' },
    @{ Row = 20; D = '14.61'; E = '  -1.60%  ' },
    @{ Row = 21; D = '5.92'; E = '  +2.11%  ' },
    @{ Row = 22; D = '438.41'; E = '  -1.90%  ' },
    @{ Row = 23; D = '0.0000197'; E = '  -2.68%  ' },
    @{ Row = 24; D = '8.99'; E = '  -5.06%  ' },
    @{ Row = 25; D = '8.25'; E = '  +4.77%  ' },
    @{ Row = 26; D = '5.64'; E = '  -0.25%  ' },
    @{ Row = 27; D = '12.55'; E = '  +6.83%  ' },
    @{ Row = 28; D = '85.51'; E = '  -2.75%  ' },
    @{ Row = 29; E = '  -0.01%  ' },
    @{ Row = 30; D = '0.183'; E = '  +9.12%  ' },
    @{ Row = 31; D = '0.249'; E = '  +5.51%  ' },
    @{ Row = 32; D = '0.123'; E = '  -15.19%  ' },
    @{ Row = 33; D = '9.12'; E = '  -2.07%  ' },
    @{ Row = 34; E = '  +0.53%  ' },
    @{ Row = 35; D = '7.88'; E = '  +1.01%  ' },
    @{ Row = 36; D = '0.158'; E = '  -10.86%  ' },
    @{ Row = 37; D = '25.72'; E = '  -1.98%  ' },
    @{ Row = 38; E = '  -5.71%  ' },
    @{ Row = 39; D = '1.89'; E = '  -2.83%  ' },
    @{ Row = 40; B = 'WhiteBITCoin'; C = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; D = '23.96'; E = '  +8.03%  ' },
    @{ Row = 41; B = 'PolygonEcosystemToken'; C = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; D = '0.440'; E = '  +0.16%  ' },
    @{ Row = 42; B = 'Fetch.AI'; C = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D = '1.28'; E = '  -2.27%  ' },
    @{ Row = 43; D = '469.18'; E = '  -4.43%  ' },
    @{ Row = 44; D = '3.24'; E = '  -5.05%  ' },
    @{ Row = 45; E = '  +0.04%  ' },
    @{ Row = 46; D = '158.88'; E = '  +0.75%  ' },
    @{ Row = 47; D = '0.695'; E = '  -1.73%  ' },
    @{ Row = 48; D = '1.84'; E = '  -3.62%  ' },
    @{ Row = 49; D = '1.31'; E = '  -3.26%  ' },
    @{ Row = 50; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '43.78'; E = '  -0.83%  ' },
    @{ Row = 51; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.34'; E = '  -1.04%  ' }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("B")) {
        $ws.Range("B$rowNum").Value = $r.B
    }
    if ($r.ContainsKey("C")) {
        $ws.Range("C$rowNum").Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        $cell = $ws.Range("D$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
        $cell.ClearFormats()
    }
    if ($r.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $r.E
    }
}
